$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 15 and 16
$ws.Range("B15").Value = 7004588
$ws.Range("B16").Value = 7004589
$ws.Range("E15").Value = "Umm Salal"
$ws.Range("E16").Value = "AlMuaidar"
$ws.Range("F15").Value = "Qatar SC Doha"
$ws.Range("F16").Value = "Al Markhiya"
$ws.Range("G15").Value = 2
$ws.Range("G16").Value = 5
$ws.Range("H15").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("I15").Value = 1
$ws.Range("I16").Value = 5
$ws.Range("J15").Value = 2
$ws.Range("J16").Value = 1
$ws.Range("K15").Value = "D"
$ws.Range("K16").Value = "H"
$ws.Range("L15").Value = 2.9
$ws.Range("L16").Value = 2.4
$ws.Range("M15").Value = 3.3
$ws.Range("M16").Value = 3.6
$ws.Range("N15").Value = 2.25
$ws.Range("N16").Value = 2.375
$ws.Range("O15").Value = 3.75
$ws.Range("O16").Value = 2.5
$ws.Range("P15").Value = 3.4
$ws.Range("P16").Value = 3.5
$ws.Range("Q15").Value = 1.909
$ws.Range("Q16").Value = 2.3
$ws.Range("R15").Value = 0.5
$ws.Range("R16").Value = 0
$ws.Range("S15").Value = 1.85
$ws.Range("S16").Value = 1.975
$ws.Range("T15").Value = 1.95
$ws.Range("T16").Value = 1.825
$ws.Range("U15").Value = 2.5
$ws.Range("U16").Value = 3
$ws.Range("V15").Value = 1.85
$ws.Range("V16").Value = 2
$ws.Range("W15").Value = 1.95
$ws.Range("W16").Value = 1.8
$ws.Range("X15").Value = -1
$ws.Range("X16").Value = 1.5
$ws.Range("Y15").Value = 2.4
$ws.Range("Y16").Value = -1
$ws.Range("Z15").Value = -1
$ws.Range("Z16").Value = -1
$ws.Range("AA15").Value = 0.8500000000000001
$ws.Range("AA16").Value = 0.9750000000000001
$ws.Range("AB15").Value = -1
$ws.Range("AB16").Value = -1
$ws.Range("AC15").Value = 0.8500000000000001
$ws.Range("AC16").Value = 1
$ws.Range("AD15").Value = -1
$ws.Range("AD16").Value = -1

# Swap rows 18 and 19
$ws.Range("B18").Value = 7004591
$ws.Range("B19").Value = 7003585
$ws.Range("E18").Value = "AlShamal SC"
$ws.Range("E19").Value = "Al Sadd"
$ws.Range("F18").Value = "AlRayyan SC"
$ws.Range("F19").Value = "AlWakrah SC"
$ws.Range("G18").Value = 3
$ws.Range("G19").Value = 0
$ws.Range("H18").Value = 4
$ws.Range("H19").Value = 0
$ws.Range("I18").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("J18").Value = 2
$ws.Range("J19").Value = 0
$ws.Range("K18").Value = "A"
$ws.Range("K19").Value = "D"
$ws.Range("L18").Value = 4.5
$ws.Range("L19").Value = 1.615
$ws.Range("M18").Value = 4.2
$ws.Range("M19").Value = 4
$ws.Range("N18").Value = 1.55
$ws.Range("N19").Value = 4.333
$ws.Range("O18").Value = 3.3
$ws.Range("O19").Value = 1.533
$ws.Range("P18").Value = 3.8
$ws.Range("P19").Value = 4.2
$ws.Range("Q18").Value = 1.85
$ws.Range("Q19").Value = 5
$ws.Range("R18").Value = 0.5
$ws.Range("R19").Value = -1
$ws.Range("S18").Value = 1.85
$ws.Range("S19").Value = 1.8
$ws.Range("T18").Value = 1.95
$ws.Range("T19").Value = 2
$ws.Range("U18").Value = 2.75
$ws.Range("U19").Value = 3.5
$ws.Range("V18").Value = 1.85
$ws.Range("V19").Value = 1.925
$ws.Range("W18").Value = 1.95
$ws.Range("W19").Value = 1.875
$ws.Range("X18").Value = -1
$ws.Range("X19").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Y19").Value = 3.2
$ws.Range("Z18").Value = 0.8500000000000001
$ws.Range("Z19").Value = -1
$ws.Range("AA18").Value = -1
$ws.Range("AA19").Value = -1
$ws.Range("AB18").Value = 0.95
$ws.Range("AB19").Value = 1
$ws.Range("AC18").Value = 0.8500000000000001
$ws.Range("AC19").Value = -1
$ws.Range("AD18").Value = -1
$ws.Range("AD19").Value = 0.875

# Swap rows 54 and 55
$ws.Range("B54").Value = 7003592
$ws.Range("B55").Value = 7004616
$ws.Range("E54").Value = "Al Sadd"
$ws.Range("E55").Value = "Al Duhail"
$ws.Range("F54").Value = "AlShamal SC"
$ws.Range("F55").Value = "Umm Salal"
$ws.Range("G54").Value = 4
$ws.Range("G55").Value = 0
$ws.Range("H54").Value = 0
$ws.Range("H55").Value = 1
$ws.Range("I54").Value = 1
$ws.Range("I55").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K54").Value = "H"
$ws.Range("K55").Value = "A"
$ws.Range("L54").Value = 1.25
$ws.Range("L55").Value = 1.285
$ws.Range("M54").Value = 5.75
$ws.Range("M55").Value = 6.5
$ws.Range("N54").Value = 8
$ws.Range("N55").Value = 6.5
$ws.Range("O54").Value = 1.25
$ws.Range("O55").Value = 1.363
$ws.Range("P54").Value = 6
$ws.Range("P55").Value = 5.75
$ws.Range("Q54").Value = 8.5
$ws.Range("Q55").Value = 5.5
$ws.Range("R54").Value = -1.75
$ws.Range("R55").Value = -1.25
$ws.Range("S54").Value = 1.825
$ws.Range("S55").Value = 1.75
$ws.Range("T54").Value = 1.975
$ws.Range("T55").Value = 1.95
$ws.Range("U54").Value = 3.75
$ws.Range("U55").Value = 3.5
$ws.Range("V54").Value = 1.975
$ws.Range("V55").Value = 1.95
$ws.Range("W54").Value = 1.825
$ws.Range("W55").Value = 1.85
$ws.Range("X54").Value = 0.25
$ws.Range("X55").Value = -1
$ws.Range("Y54").Value = -1
$ws.Range("Y55").Value = -1
$ws.Range("Z54").Value = -1
$ws.Range("Z55").Value = 4.5
$ws.Range("AA54").Value = 0.825
$ws.Range("AA55").Value = -1
$ws.Range("AB54").Value = -1
$ws.Range("AB55").Value = 0.95
$ws.Range("AC54").Value = 0.4875
$ws.Range("AC55").Value = -1
$ws.Range("AD54").Value = -0.5
$ws.Range("AD55").Value = 0.8500000000000001

# Swap rows 60 and 61
$ws.Range("B60").Value = 7003485
$ws.Range("B61").Value = 7004620
$ws.Range("E60").Value = "AlWakrah SC"
$ws.Range("E61").Value = "Umm Salal"
$ws.Range("F60").Value = "AlAhli Doha"
$ws.Range("F61").Value = "Al Markhiya"
$ws.Range("G60").Value = 3
$ws.Range("G61").Value = 1
$ws.Range("H60").Value = 2
$ws.Range("H61").Value = 1
$ws.Range("I60").Value = 2
$ws.Range("I61").Value = 1
$ws.Range("J60").Value = 1
$ws.Range("J61").Value = 1
$ws.Range("K60").Value = "H"
$ws.Range("K61").Value = "D"
$ws.Range("L60").Value = 1.533
$ws.Range("L61").Value = 1.727
$ws.Range("M60").Value = 4.2
$ws.Range("M61").Value = 3.75
$ws.Range("N60").Value = 5
$ws.Range("N61").Value = 4
$ws.Range("O60").Value = 1.615
$ws.Range("O61").Value = 1.833
$ws.Range("P60").Value = 3.75
$ws.Range("P61").Value = 3.6
$ws.Range("Q60").Value = 4.75
$ws.Range("Q61").Value = 3.75
$ws.Range("R60").Value = -0.75
$ws.Range("R61").Value = -0.5
$ws.Range("S60").Value = 1.8
$ws.Range("S61").Value = 1.9
$ws.Range("T60").Value = 2
$ws.Range("T61").Value = 1.9
$ws.Range("U60").Value = 3
$ws.Range("U61").Value = 3
$ws.Range("V60").Value = 1.975
$ws.Range("V61").Value = 2
$ws.Range("W60").Value = 1.825
$ws.Range("W61").Value = 1.8
$ws.Range("X60").Value = 0.615
$ws.Range("X61").Value = -1
$ws.Range("Y60").Value = -1
$ws.Range("Y61").Value = 2.6
$ws.Range("Z60").Value = -1
$ws.Range("Z61").Value = -1
$ws.Range("AA60").Value = 0.4
$ws.Range("AA61").Value = -1
$ws.Range("AB60").Value = -0.5
$ws.Range("AB61").Value = 0.8999999999999999
$ws.Range("AC60").Value = 0.9750000000000001
$ws.Range("AC61").Value = -1
$ws.Range("AD60").Value = -1
$ws.Range("AD61").Value = 0.8

# Swap rows 87 and 88
$ws.Range("B87").Value = 7840806
$ws.Range("B88").Value = 7840807
$ws.Range("E87").Value = "Umm Salal"
$ws.Range("E88").Value = "AlMuaidar"
$ws.Range("F87").Value = "AlShamal SC"
$ws.Range("F88").Value = "AlAhli Doha"
$ws.Range("G87").Value = 3
$ws.Range("G88").Value = 1
$ws.Range("H87").Value = 4
$ws.Range("H88").Value = 1
$ws.Range("I87").Value = 0
$ws.Range("I88").Value = 1
$ws.Range("J87").Value = 2
$ws.Range("J88").Value = 0
$ws.Range("K87").Value = "A"
$ws.Range("K88").Value = "D"
$ws.Range("L87").Value = 2
$ws.Range("L88").Value = 2.875
$ws.Range("M87").Value = 3.6
$ws.Range("M88").Value = 4
$ws.Range("N87").Value = 3.2
$ws.Range("N88").Value = 2
$ws.Range("O87").Value = 2.2
$ws.Range("O88").Value = 1.85
$ws.Range("P87").Value = 3.5
$ws.Range("P88").Value = 3.8
$ws.Range("Q87").Value = 2.875
$ws.Range("Q88").Value = 3.4
$ws.Range("R87").Value = -0.25
$ws.Range("R88").Value = -0.5
$ws.Range("S87").Value = 1.95
$ws.Range("S88").Value = 1.875
$ws.Range("T87").Value = 1.85
$ws.Range("T88").Value = 1.925
$ws.Range("U87").Value = 2.75
$ws.Range("U88").Value = 3.25
$ws.Range("V87").Value = 1.825
$ws.Range("V88").Value = 1.95
$ws.Range("W87").Value = 1.975
$ws.Range("W88").Value = 1.75
$ws.Range("X87").Value = -1
$ws.Range("X88").Value = -1
$ws.Range("Y87").Value = -1
$ws.Range("Y88").Value = 2.8
$ws.Range("Z87").Value = 1.875
$ws.Range("Z88").Value = -1
$ws.Range("AA87").Value = -1
$ws.Range("AA88").Value = -1
$ws.Range("AB87").Value = 0.8500000000000001
$ws.Range("AB88").Value = 0.925
$ws.Range("AC87").Value = 0.825
$ws.Range("AC88").Value = -1
$ws.Range("AD87").Value = -1
$ws.Range("AD88").Value = 0.75

# Swap rows 94 and 95
$ws.Range("B94").Value = 7936025
$ws.Range("B95").Value = 7936026
$ws.Range("E94").Value = "AlMuaidar"
$ws.Range("E95").Value = "AlShamal SC"
$ws.Range("F94").Value = "Al Gharafa"
$ws.Range("F95").Value = "Al Markhiya"
$ws.Range("G94").Value = 0
$ws.Range("G95").Value = 2
$ws.Range("H94").Value = 2
$ws.Range("H95").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K94").Value = "A"
$ws.Range("K95").Value = "H"
$ws.Range("L94").Value = 4.75
$ws.Range("L95").Value = 2.375
$ws.Range("M94").Value = 4.333
$ws.Range("M95").Value = 3.3
$ws.Range("N94").Value = 1.55
$ws.Range("N95").Value = 2.8
$ws.Range("O94").Value = 4.333
$ws.Range("O95").Value = 2.05
$ws.Range("P94").Value = 3.8
$ws.Range("P95").Value = 3.4
$ws.Range("Q94").Value = 1.666
$ws.Range("Q95").Value = 3.3
$ws.Range("R94").Value = 0.75
$ws.Range("R95").Value = -0.25
$ws.Range("S94").Value = 1.95
$ws.Range("S95").Value = 1.8
$ws.Range("T94").Value = 1.85
$ws.Range("T95").Value = 2
$ws.Range("U94").Value = 3
$ws.Range("U95").Value = 2.75
$ws.Range("V94").Value = 1.825
$ws.Range("V95").Value = 1.9
$ws.Range("W94").Value = 1.975
$ws.Range("W95").Value = 1.9
$ws.Range("X94").Value = -1
$ws.Range("X95").Value = 1.05
$ws.Range("Y94").Value = -1
$ws.Range("Y95").Value = -1
$ws.Range("Z94").Value = 0.6659999999999999
$ws.Range("Z95").Value = -1
$ws.Range("AA94").Value = -1
$ws.Range("AA95").Value = 0.8
$ws.Range("AB94").Value = 0.8500000000000001
$ws.Range("AB95").Value = -1
$ws.Range("AC94").Value = -1
$ws.Range("AC95").Value = -1
$ws.Range("AD94").Value = 0.9750000000000001
$ws.Range("AD95").Value = 0.8999999999999999

